$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark currently sitting on the
#    "Jelszo torlese" heading paragraph (it will be re-created further
#    down, at the end of the new content being added below).
# ------------------------------------------------------------------
try {
    $oldBookmark = $d.Bookmarks.Item("_GoBack")
    $oldBookmark.Delete()
} catch {
    # no pre-existing _GoBack bookmark - nothing to remove
}

# ------------------------------------------------------------------
# 2) Locate the paragraph that ends the "Jelszo torles" use case
#    ("Az alkalmazas kitorli a jelszo adatait") - the new "Jelszo
#    modositasa" use case is inserted right after it, before "Mappa
#    hozzaadas".
# ------------------------------------------------------------------
$anchor = "kit" + [char]0x00F6 + "rli a jelsz" + [char]0x00F3 + " adatait"

$findRange = $d.Content
$found = $findRange.Find.Execute($anchor, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate anchor paragraph for insertion"
}

$insertPoint = $d.Range($findRange.End, $findRange.End)

# ------------------------------------------------------------------
# 3) Build the OOXML for the six new paragraphs (mirrors the other
#    use-case blocks: a bold heading + Aktorok + Foforgatokonyv +
#    three numbered steps). A trailing empty <w:p/> is appended
#    because InsertXML merges the *last* fragment paragraph into
#    whatever paragraph follows the insertion point - the empty
#    paragraph absorbs that merge so "Mappa hozzaadas" stays intact,
#    and is then deleted again afterwards.
# ------------------------------------------------------------------
$newBlock = '<w:p><w:pPr><w:pStyle w:val="Listaszerbekezds"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Jelszó módosítása</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listaszerbekezds"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Aktorok: Felhasználó, alkalmazás</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listaszerbekezds"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Főforgatókönyv:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listaszerbekezds"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>A felhasználó kiválaszt egy mappa jelszavait és rákattint</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listaszerbekezds"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>A felhasználó a táblázatban átírja az egyik jelszó adatait</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listaszerbekezds"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Az alkalmazás elmenti a jelszóhoz tartozó új adatokat</w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/></w:p><w:p/>'

$payload = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    $newBlock + `
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$countBefore = $d.Paragraphs.Count
$insertPoint.InsertXML($payload)
$countAfter = $d.Paragraphs.Count

# ------------------------------------------------------------------
# 4) Delete the trailing placeholder paragraph that absorbed the
#    merge described above - it is the last of the newly-added
#    paragraphs.
# ------------------------------------------------------------------
$addedCount = $countAfter - $countBefore
if ($addedCount -gt 0) {
    $strayIndex = $countAfter - ($addedCount - 6)
    $stray = $d.Paragraphs.Item($strayIndex)
    if ($stray.Range.Text -eq [char]13) {
        $stray.Range.Delete()
    }
}
